# Update cryptocurrency price (D) and volume-change (E) columns
# for rows 2-51 to reflect the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.817.19"
$ws.Range("E2").Value = "  -1.29%  "

# Row 3
$ws.Range("D3").Value = "1.744.44"
$ws.Range("E3").Value = "  -2.03%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3877"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.44%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3377"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.76%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.82%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07156"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.58%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.55%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.072"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.59%  "

# Row 15
$ws.Range("D15").Value = "1.740.42"
$ws.Range("E15").Value = "  -2.62%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.946"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.21%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001048"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.97%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06594"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "79.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.80%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.29%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.147"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.56%  "

# Row 23
$ws.Range("D23").Value = "27.765.87"
$ws.Range("E23").Value = "  -1.68%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.03%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.382"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.54%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.272"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.97%  "

# Row 29
$ws.Range("D29").Value = "1.939.45"
$ws.Range("E29").Value = "  -2.49%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.261"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.30%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.50%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.076"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.80%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.749"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.54%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08705"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.98%  "

# Row 36
$ws.Range("E36").Value = "  +0.89%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02259"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.34%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.079"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.79%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06078"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.30%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6372"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.16%  "

# Row 41
$ws.Range("E41").Value = "  -4.79%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.187"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.33%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.25%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.837"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.85%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.55%  "

# Row 46
$ws.Range("E46").Value = "  -1.57%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5903"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.34%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.76%  "

# Row 49
$ws.Range("E49").Value = "  -6.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06934"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.140"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.13%  "
